$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the "✓ "/"✗ " glyph prefixes from the Is Active column text
#    ("✓ Active" -> "Active", "✗ Inactive" -> "Inactive").
# ------------------------------------------------------------------
$ws.Range("F2").Value2 = "Active"
$ws.Range("F3").Value2 = "Active"
$ws.Range("F4").Value2 = "Inactive"
$ws.Range("F5").Value2 = "Active"

# ------------------------------------------------------------------
# 2. Fix wrapping: the data rows should wrap text and left-align
#    instead of "general" horizontal alignment with no wrap.
# ------------------------------------------------------------------
$data = $ws.Range("A2:G5")
$data.HorizontalAlignment = -4131   # xlLeft
$data.WrapText = $true

# ------------------------------------------------------------------
# 3. Fix heights: widen every column by 2 characters so the wrapped
#    text has room to breathe.
# ------------------------------------------------------------------
for ($i = 1; $i -le 7; $i++) {
    $col = $ws.Columns.Item($i)
    $col.ColumnWidth = $col.ColumnWidth + 2
}
